$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and G columns (data rows 2-51) are treated as text so that
# numeric-looking strings (with significant trailing/leading zeros) are
# preserved exactly as text rather than being coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2,4).Value = "290.18"
$ws.Cells.Item(2,7).Value = "19"

# Row 3
$ws.Cells.Item(3,4).Value = "21.35"
$ws.Cells.Item(3,7).Value = "19"

# Row 4
$ws.Cells.Item(4,4).Value = "6.465"
$ws.Cells.Item(4,7).Value = "19"

# Row 5
$ws.Cells.Item(5,4).Value = "0.06402"
$ws.Cells.Item(5,7).Value = "19"

# Row 6
$ws.Cells.Item(6,2).Value = "FTXToken"
$ws.Cells.Item(6,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(6,4).Value = "1.636"
$ws.Cells.Item(6,5).Value = "5FTXTokenFTT"
$ws.Cells.Item(6,7).Value = "19"

# Row 7
$ws.Cells.Item(7,2).Value = "GateToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(7,4).Value = "3.615"
$ws.Cells.Item(7,5).Value = "6GateTokenGT"
$ws.Cells.Item(7,7).Value = "19"

# Row 8
$ws.Cells.Item(8,4).Value = "6.621"
$ws.Cells.Item(8,7).Value = "19"

# Row 9
$ws.Cells.Item(9,4).Value = "0.8292"
$ws.Cells.Item(9,7).Value = "19"

# Row 10
$ws.Cells.Item(10,4).Value = "0.01439"
$ws.Cells.Item(10,7).Value = "19"

# Row 11
$ws.Cells.Item(11,4).Value = "0.1705"
$ws.Cells.Item(11,7).Value = "19"

# Row 12
$ws.Cells.Item(12,4).Value = "0.08747"
$ws.Cells.Item(12,7).Value = "19"

# Row 13
$ws.Cells.Item(13,4).Value = "0.03662"
$ws.Cells.Item(13,7).Value = "19"

# Row 14
$ws.Cells.Item(14,4).Value = "0.03208"
$ws.Cells.Item(14,7).Value = "19"

# Row 15
$ws.Cells.Item(15,4).Value = "0.09209"
$ws.Cells.Item(15,7).Value = "19"

# Row 16
$ws.Cells.Item(16,4).Value = "3.705"
$ws.Cells.Item(16,7).Value = "19"

# Row 17
$ws.Cells.Item(17,4).Value = "0.001660"
$ws.Cells.Item(17,7).Value = "19"

# Row 18
$ws.Cells.Item(18,4).Value = "0.04753"
$ws.Cells.Item(18,7).Value = "19"

# Row 19
$ws.Cells.Item(19,4).Value = "0.006206"
$ws.Cells.Item(19,7).Value = "19"

# Row 20
$ws.Cells.Item(20,4).Value = "0.005745"
$ws.Cells.Item(20,7).Value = "19"

# Row 21
$ws.Cells.Item(21,4).Value = "0.001073"
$ws.Cells.Item(21,7).Value = "19"

# Row 22
$ws.Cells.Item(22,4).Value = "0.0001602"
$ws.Cells.Item(22,7).Value = "19"

# Row 23
$ws.Cells.Item(23,4).Value = "3.761"
$ws.Cells.Item(23,7).Value = "19"

# Row 24
$ws.Cells.Item(24,4).Value = "2.406"
$ws.Cells.Item(24,7).Value = "19"

# Row 25
$ws.Cells.Item(25,4).Value = "0.3349"
$ws.Cells.Item(25,7).Value = "19"

# Row 26
$ws.Cells.Item(26,4).Value = "0.1261"
$ws.Cells.Item(26,7).Value = "19"

# Row 27
$ws.Cells.Item(27,7).Value = "19"

# Row 28
$ws.Cells.Item(28,7).Value = "19"

# Row 29
$ws.Cells.Item(29,7).Value = "19"

# Row 30
$ws.Cells.Item(30,7).Value = "19"

# Row 31
$ws.Cells.Item(31,7).Value = "19"

# Row 32
$ws.Cells.Item(32,7).Value = "19"

# Row 33
$ws.Cells.Item(33,7).Value = "19"

# Row 34
$ws.Cells.Item(34,7).Value = "19"

# Row 35
$ws.Cells.Item(35,7).Value = "19"

# Row 36
$ws.Cells.Item(36,7).Value = "19"

# Row 37
$ws.Cells.Item(37,7).Value = "19"

# Row 38
$ws.Cells.Item(38,7).Value = "19"

# Row 39
$ws.Cells.Item(39,7).Value = "19"

# Row 40
$ws.Cells.Item(40,4).Value = "0.04900"
$ws.Cells.Item(40,7).Value = "19"

# Row 41
$ws.Cells.Item(41,4).Value = "0.007159"
$ws.Cells.Item(41,7).Value = "19"

# Row 42
$ws.Cells.Item(42,4).Value = "0.004505"
$ws.Cells.Item(42,7).Value = "19"

# Row 43
$ws.Cells.Item(43,4).Value = "0.1127"
$ws.Cells.Item(43,7).Value = "19"

# Row 44
$ws.Cells.Item(44,4).Value = "0.01149"
$ws.Cells.Item(44,7).Value = "19"

# Row 45
$ws.Cells.Item(45,4).Value = "0.00006970"
$ws.Cells.Item(45,7).Value = "19"

# Row 46
$ws.Cells.Item(46,4).Value = "0.00000000751"
$ws.Cells.Item(46,7).Value = "19"

# Row 47
$ws.Cells.Item(47,4).Value = "0.8007"
$ws.Cells.Item(47,7).Value = "19"

# Row 48
$ws.Cells.Item(48,4).Value = "0.007319"
$ws.Cells.Item(48,7).Value = "19"

# Row 49
$ws.Cells.Item(49,4).Value = "0.00001902"
$ws.Cells.Item(49,7).Value = "19"

# Row 50
$ws.Cells.Item(50,7).Value = "19"

# Row 51
$ws.Cells.Item(51,7).Value = "19"
